$wb = $excel.ActiveWorkbook

# New file handed back: 27f9817c-2f0d-4784-ad72-521c8ab9f3b9.md
# It is inserted as row 3 in every sheet (Overview, zh-cn, de-de),
# pushing the previously-row-3 entry (e53859ba-...) down to row 4.

$newBase   = "27f9817c-2f0d-4784-ad72-521c8ab9f3b9"
$oldBase   = "e53859ba-c3c5-40b4-be9a-ee0450456606"
$otherBase = "054b27fe-8ff8-4fc9-a088-9197af2b1b7f"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "$newBase.md"
$ws.Range("B3").Value = "e2e\$newBase.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-09-04 08:48:09"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G4"))

# Hyperlinks got desynced by the row insert (their anchor stayed put
# while the underlying data moved down) - rebuild all of them fresh.
$ws.Range("A1").Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/382ef41d03b767a9236653e29c53aaf8252d314c/e2e/$otherBase.md", [Type]::Missing, [Type]::Missing, "e2e\$otherBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a660a2ba04b449ecf904d8533d7db0df197a93f1/e2e/$newBase.md", [Type]::Missing, [Type]::Missing, "e2e\$newBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ecde17806d6266c93a507b97622aff58f57b7ef/e2e/$oldBase.md", [Type]::Missing, [Type]::Missing, "e2e\$oldBase.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "$newBase.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newBase.a660a2ba04b449ecf904d8533d7db0df197a93f1.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-04 08:47:59"
$ws.Range("I3").Value = "$newBase.md"
$ws.Range("J3").Value = "$newBase.a660a2ba04b449ecf904d8533d7db0df197a93f1.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-04 08:48:27"
$ws.Range("L3").Value = "False"
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = "False"
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = "False"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P4"))

$ws.Range("A1").Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/382ef41d03b767a9236653e29c53aaf8252d314c/e2e/$otherBase.md", [Type]::Missing, [Type]::Missing, "$otherBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5eba0a01eeee2683353807709a0387ff40f7bfb9/e2e/$otherBase.md", [Type]::Missing, [Type]::Missing, "$otherBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a660a2ba04b449ecf904d8533d7db0df197a93f1/e2e/$newBase.md", [Type]::Missing, [Type]::Missing, "$newBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a660a2ba04b449ecf904d8533d7db0df197a93f1/e2e/$newBase.md", [Type]::Missing, [Type]::Missing, "$newBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ecde17806d6266c93a507b97622aff58f57b7ef/e2e/$oldBase.md", [Type]::Missing, [Type]::Missing, "$oldBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4ddb87a16996585d70af82a13add249cd6dc70a9/e2e/$oldBase.md", [Type]::Missing, [Type]::Missing, "$oldBase.md")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "$newBase.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newBase.a660a2ba04b449ecf904d8533d7db0df197a93f1.de-de.xlf"
$ws.Range("H3").Value = "2016-09-04 08:48:09"
$ws.Range("I3").Value = "$newBase.md"
$ws.Range("J3").Value = "$newBase.a660a2ba04b449ecf904d8533d7db0df197a93f1.de-de.xlf"
$ws.Range("K3").Value = "2016-09-04 08:48:34"
$ws.Range("L3").Value = "False"
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = "False"
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = "False"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P4"))

$ws.Range("A1").Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/382ef41d03b767a9236653e29c53aaf8252d314c/e2e/$otherBase.md", [Type]::Missing, [Type]::Missing, "$otherBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a40ba5e4024fe8990f45716f8d437c6ae3520d7c/e2e/$otherBase.md", [Type]::Missing, [Type]::Missing, "$otherBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a660a2ba04b449ecf904d8533d7db0df197a93f1/e2e/$newBase.md", [Type]::Missing, [Type]::Missing, "$newBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a660a2ba04b449ecf904d8533d7db0df197a93f1/e2e/$newBase.md", [Type]::Missing, [Type]::Missing, "$newBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ecde17806d6266c93a507b97622aff58f57b7ef/e2e/$oldBase.md", [Type]::Missing, [Type]::Missing, "$oldBase.md")
$null = $ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/61ec165686a4665caa3203b81540300780876901/e2e/$oldBase.md", [Type]::Missing, [Type]::Missing, "$oldBase.md")

$wb.Worksheets.Item("Overview").Activate()
